# User Story 08 - Birth Before Marriage
# Mark US04 (row 5) and US08 (row 7) of the Sprint1 sheet as Done,
# filling in actual size/time + completion date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint1")

# Row 5 -> US04 "Marriage before divorce"
$ws.Range("D5").Value = "Done"
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = 40
$ws.Range("I5").NumberFormat = $ws.Range("I2").NumberFormat
$ws.Range("I5").Value = (Get-Date -Year 2021 -Month 10 -Day 3).Date

# Row 7 -> US08 "Birth before marriage of parents"
$ws.Range("D7").Value = "Done"
$ws.Range("G7").Value = 35
$ws.Range("H7").Value = 25
$ws.Range("I7").NumberFormat = $ws.Range("I2").NumberFormat
$ws.Range("I7").Value = (Get-Date -Year 2021 -Month 10 -Day 4).Date

# Make Sprint1 the active/selected sheet, with I7 as the active cell
$ws.Activate()
$ws.Range("I7").Select()
